# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 164

    $ws.Range("F3").Value = 7374
    $ws.Range("G3").Value = 54

    if ($name -eq "展览") {
        $ws.Range("F4").Value = 7013
    } else {
        $ws.Range("F4").Value = 7018
    }

    $ws.Range("F5").Value = 88
    $ws.Range("F6").Value = 180
    $ws.Range("F7").Value = 35
    $ws.Range("F9").Value = 114
    $ws.Range("F10").Value = 91
    $ws.Range("F11").Value = 122
    $ws.Range("F12").Value = 209
    $ws.Range("F13").Value = 82
    $ws.Range("F14").Value = 658
    $ws.Range("F15").Value = 458
    $ws.Range("F17").Value = 20
    $ws.Range("F20").Value = 67
}
